$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 124, shifting rows 124:143 down to 125:144.
$ws.Rows.Item(124).Insert()

# Populate the newly inserted row 124 with the new record's data.
$ws.Cells.Item(124, 1).Value = 10
$ws.Cells.Item(124, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(124, 3).Value = "La Araucanía"
$ws.Cells.Item(124, 4).Value = 44694
$ws.Cells.Item(124, 5).Value = 9
$ws.Cells.Item(124, 6).Value = 100112012
$ws.Cells.Item(124, 7).Value = "Espinaca"
$ws.Cells.Item(124, 8).Value = "Sin especificar"
$ws.Cells.Item(124, 9).Value = "Primera"
$ws.Cells.Item(124, 10).Value = 85
$ws.Cells.Item(124, 11).Value = 9000
$ws.Cells.Item(124, 12).Value = 9000
$ws.Cells.Item(124, 13).Value = 9000
$ws.Cells.Item(124, 14).Value = "$/docena de atados"
$ws.Cells.Item(124, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(124, 16).Value = 3000
$ws.Cells.Item(124, 17).Value = 3
$ws.Cells.Item(124, 18).Value = "Hortaliza"
